$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.027979016304016
$ws.Range("B1").Value = 2.297914505004883
$ws.Range("C1").Value = 4.657309532165527
$ws.Range("D1").Value = 1.367353916168213
$ws.Range("E1").Value = 1.269149780273438
